# Update odds values per diff for Jogos_da_Semana_FlashScore_2024-10-03.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 5.25
$ws.Range("J2").Value = 2.5
$ws.Range("K2").Value = 1.91
$ws.Range("L2").Value = 6.5
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 2.38
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.48
$ws.Range("S2").Value = 1.62
$ws.Range("T2").Value = 2.2
$ws.Range("W2").Value = 4.75
$ws.Range("X2").Value = 6.5
$ws.Range("Z2").Value = 13
$ws.Range("AA2").Value = 19
$ws.Range("AC2").Value = 6
$ws.Range("AD2").Value = 7
$ws.Range("AF2").Value = 101
$ws.Range("AG2").Value = 9.5
$ws.Range("AH2").Value = 26
$ws.Range("AJ2").Value = 67
$ws.Range("AN2").Value = 3.4
$ws.Range("AO2").Value = 10
$ws.Range("AP2").Value = 29
$ws.Range("AQ2").Value = 34
$ws.Range("AR2").Value = 67
$ws.Range("AS2").Value = 301
$ws.Range("AT2").Value = 2.2
$ws.Range("AW2").Value = 7
# Row 3
$ws.Range("G3").Value = 1.6
$ws.Range("H3").Value = 3.9
$ws.Range("K3").Value = 2.1
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 2.15
$ws.Range("R3").Value = 1.67
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 6.5
$ws.Range("Y3").Value = 9
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 8.5
$ws.Range("AF3").Value = 81
$ws.Range("AG3").Value = 12
$ws.Range("AH3").Value = 26
$ws.Range("AK3").Value = 51
$ws.Range("AO3").Value = 8.5
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 2.63
$ws.Range("AU3").Value = 9.5
# Row 4
$ws.Range("G4").Value = 1.67
$ws.Range("H4").Value = 3.2
$ws.Range("J4").Value = 2.4
$ws.Range("K4").Value = 1.91
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 2.2
$ws.Range("Q4").Value = 2.88
$ws.Range("R4").Value = 1.4
$ws.Range("S4").Value = 1.62
$ws.Range("T4").Value = 2.2
$ws.Range("U4").Value = 2.63
$ws.Range("V4").Value = 1.44
$ws.Range("W4").Value = 4.5
$ws.Range("Y4").Value = 10
$ws.Range("Z4").Value = 12
$ws.Range("AC4").Value = 5.5
$ws.Range("AF4").Value = 126
$ws.Range("AI4").Value = 23
$ws.Range("AK4").Value = 67
$ws.Range("AL4").Value = 81
$ws.Range("AO4").Value = 9.5
$ws.Range("AR4").Value = 81
$ws.Range("AS4").Value = 351
$ws.Range("AT4").Value = 2.2
# Row 5
$ws.Range("Q5").Value = 1.9
$ws.Range("R5").Value = 1.9
$ws.Range("AM5").Value = 1000
# Row 6
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 1.53
$ws.Range("J6").Value = 4.75
$ws.Range("K6").Value = 2.5
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 11
$ws.Range("O6").Value = 1.14
$ws.Range("P6").Value = 5
$ws.Range("Q6").Value = 1.5
$ws.Range("R6").Value = 2.5
$ws.Range("S6").Value = 1.25
$ws.Range("T6").Value = 3.75
$ws.Range("U6").Value = 1.57
$ws.Range("V6").Value = 2.25
$ws.Range("W6").Value = 21
$ws.Range("Y6").Value = 17
$ws.Range("AA6").Value = 34
$ws.Range("AB6").Value = 34
$ws.Range("AC6").Value = 19
$ws.Range("AD6").Value = 8.5
$ws.Range("AE6").Value = 13
$ws.Range("AG6").Value = 10
$ws.Range("AH6").Value = 9.5
$ws.Range("AI6").Value = 9
$ws.Range("AJ6").Value = 12
$ws.Range("AL6").Value = 19
$ws.Range("AM6").Value = 126
$ws.Range("AN6").Value = 7
$ws.Range("AO6").Value = 23
$ws.Range("AP6").Value = 26
$ws.Range("AQ6").Value = 67
$ws.Range("AR6").Value = 81
$ws.Range("AS6").Value = 126
$ws.Range("AT6").Value = 3.75
$ws.Range("AU6").Value = 7.5
$ws.Range("AV6").Value = 41
$ws.Range("AW6").Value = 4
$ws.Range("AY6").Value = 15
$ws.Range("BA6").Value = 34
$ws.Range("BB6").Value = 81
# Row 8
$ws.Range("G8").Value = 2.8
$ws.Range("I8").Value = 2.35
$ws.Range("J8").Value = 3.25
$ws.Range("L8").Value = 3
$ws.Range("X8").Value = 15
$ws.Range("AQ8").Value = 41
$ws.Range("AX8").Value = 13
$ws.Range("BB8").Value = 126
